# Remove "img" prefix and ".png" extension from the ImageName column (B),
# e.g. "AuraSkypoolDXB-img1.png" -> "AuraSkypoolDXB-1", for rows 2..213.
# Also fixes a duplicate HolderAddress on row 201 (was accidentally the same
# as row 200) by giving it its own unique address.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 213; $row++) {
    $n = $row - 1
    $ws.Cells.Item($row, 2).Value = "AuraSkypoolDXB-$n"
}

$ws.Cells.Item(201, 1).Value = "0x3a4fdd7a51d388218a33b559a8c1a67f24791e6c"
